$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G (K) values for rows 2-6 as per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
